$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatted style from A6 (bold, border, centered) down to A7:A9
# so the new index cells match the existing index-column style (s="1").
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7:A9").PasteSpecial(-4122) | Out-Null

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "The result of “ab”+”c”*2 is ()"
$ws.Range("C7").Value = "abc2"
$ws.Range("D7").Value = "abcabc"
$ws.Range("E7").Value = "abcc"
$ws.Range("F7").Value = "ababcc"
$ws.Range("G7").Value = "C"

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Which of the following options is not correct about the following code  ()
fname = input(`"请输入要写入的文件: `")
fo = open(fname, `"w+`")
ls = [`"清明时节雨纷纷，`",`"路上行人欲断魂，`",`"借问酒家何处有？`",`"牧童遥指杏花村。`"]
fo.writelines(ls)
fo.seek(0)
for line in fo:
    print(line)
fo.close()"
$ws.Range("C8").Value = "fo. seek (0) can be omitted，the output is unchanged."
$ws.Range("D8").Value = "fo. writelines (ls) writes the ls list whose elements are all strings to a file"
$ws.Range("E8").Value = "The main function of the code is to write a list type to the file and print out the result"
$ws.Range("F8").Value = "When executing the code, enter `"Qingming.txt`" from the keyboard, and Qingming.txt is created"
$ws.Range("G8").Value = "A"

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "关于Python程序中与“缩进”有关的说法中，以下选项中正确的是"
$ws.Range("C9").Value = "缩进统一为4个空格"
$ws.Range("D9").Value = "缩进可以用在任何语句之后，表示语句间的包含关系"
$ws.Range("E9").Value = "缩进在程序中长度统一且强制使用"
$ws.Range("F9").Value = "缩进是非强制性的，仅为了提高代码可读性"
$ws.Range("G9").Value = "C"

